# Update the "想去人数" (F column) figures for both the "展览" and
# "全部类型" worksheets, which contain duplicate rows of the same events.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 1182
    3  = 599
    6  = 174
    7  = 62
    10 = 5505
    11 = 4896
    16 = 201
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
